# "customer 기능 1차 완료"
# Update the "Customer RestApi" sheet (4th sheet):
#  - C6: "구매 요청" -> "상품 주문" (request body changed from purchase-request to order)
#  - Append two new API rows (회원 탈퇴 / DELETE / /customer, and 입금 / PUT / /customer/balance / money)
#  - Append one trailing blank template row, matching existing row style
#  - Leave selection on I18 (matches final saved cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# 1. Extend the formatting of the last data row (34) down through the two
#    new rows and the trailing blank row (35-37) before filling in values.
$ws.Range("B34:H34").Copy()
$ws.Range("B35:H37").PasteSpecial(-4122)

# 2. New row 35 - 회원 탈퇴 (account deletion)
$ws.Range("B35").Value = 33
$ws.Range("D35").Value = "DELETE"
$ws.Range("C35").Value = "회원 탈퇴"
$ws.Range("E35").Value = "/customer"

# 3. New row 36 - 입금 (deposit / balance top-up)
$ws.Range("B36").Value = 34
$ws.Range("C36").Value = "입금"
$ws.Range("D36").Value = "PUT"
$ws.Range("E36").Value = "/customer/balance"

# 4. Fix the request-body cell text for the "상품 주문" (order) row
$ws.Range("C6").Value = "상품 주문"

# 5. Finish row 36
$ws.Range("F36").Value = "money"

# 6. New trailing blank row 37 (template row, numbered but otherwise empty)
$ws.Range("B37").Value = 35

# 7. Match the saved selection / cursor state
$ws.Activate() | Out-Null
$ws.Range("I18").Select() | Out-Null
